# Natmi following Dr Hou advice
# Update Igf2-Igf2r LR-pair stats: Ligand/Receptor-expressing cell counts (3 cells now)
# and the resulting derived expression/specificity metrics for rows 2-10.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 3.228519
$ws.Range("H2").Value = 9.685557
$ws.Range("I2").Value = 0.0641296566303666
$ws.Range("J2").Value = 0.0641296566303666
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 31.94073466666667
$ws.Range("N2").Value = 95.822204
$ws.Range("O2").Value = 0.3365562672414605
$ws.Range("P2").Value = 0.3365562672414606
$ws.Range("Q2").Value = 103.121268745292
$ws.Range("R2").Value = 928.0914187076279
$ws.Range("S2").Value = 0.02158323785499276
$ws.Range("T2").Value = 0.02158323785499277

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 3.228519
$ws.Range("H3").Value = 9.685557
$ws.Range("I3").Value = 0.0641296566303666
$ws.Range("J3").Value = 0.0641296566303666
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 38.58528233333333
$ws.Range("N3").Value = 115.755847
$ws.Range("O3").Value = 0.406569189096231
$ws.Range("P3").Value = 0.406569189096231
$ws.Range("Q3").Value = 124.573317133531
$ws.Range("R3").Value = 1121.159854201779
$ws.Range("S3").Value = 0.02607314249322788
$ws.Range("T3").Value = 0.02607314249322788

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 3.228519
$ws.Range("H4").Value = 9.685557
$ws.Range("I4").Value = 0.0641296566303666
$ws.Range("J4").Value = 0.0641296566303666
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 24.37857333333333
$ws.Range("N4").Value = 73.13571999999999
$ws.Range("O4").Value = 0.2568745436623085
$ws.Range("P4").Value = 0.2568745436623085
$ws.Range("Q4").Value = 78.70668719956
$ws.Range("R4").Value = 708.3601847960399
$ws.Range("S4").Value = 0.01647327628214596
$ws.Range("T4").Value = 0.01647327628214596

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 9.281951000000001
$ws.Range("H5").Value = 27.845853
$ws.Range("I5").Value = 0.1843719459262553
$ws.Range("J5").Value = 0.1843719459262553
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 31.94073466666667
$ws.Range("N5").Value = 95.822204
$ws.Range("O5").Value = 0.3365562672414605
$ws.Range("P5").Value = 0.3365562672414606
$ws.Range("Q5").Value = 296.4723340800014
$ws.Range("R5").Value = 2668.251006720012
$ws.Range("S5").Value = 0.0620515339049849
$ws.Range("T5").Value = 0.06205153390498491

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 9.281951000000001
$ws.Range("H6").Value = 27.845853
$ws.Range("I6").Value = 0.1843719459262553
$ws.Range("J6").Value = 0.1843719459262553
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 38.58528233333333
$ws.Range("N6").Value = 115.755847
$ws.Range("O6").Value = 0.406569189096231
$ws.Range("P6").Value = 0.406569189096231
$ws.Range("Q6").Value = 358.1466999391657
$ws.Range("R6").Value = 3223.320299452491
$ws.Range("S6").Value = 0.07495995254733177
$ws.Range("T6").Value = 0.07495995254733179

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 9.281951000000001
$ws.Range("H7").Value = 27.845853
$ws.Range("I7").Value = 0.1843719459262553
$ws.Range("J7").Value = 0.1843719459262553
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 24.37857333333333
$ws.Range("N7").Value = 73.13571999999999
$ws.Range("O7").Value = 0.2568745436623085
$ws.Range("P7").Value = 0.2568745436623085
$ws.Range("Q7").Value = 226.2807231299067
$ws.Range("R7").Value = 2036.52650816916
$ws.Range("S7").Value = 0.04736045947393866
$ws.Range("T7").Value = 0.04736045947393866

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 37.833149
$ws.Range("H8").Value = 113.499447
$ws.Range("I8").Value = 0.751498397443378
$ws.Range("J8").Value = 0.7514983974433781
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 31.94073466666667
$ws.Range("N8").Value = 95.822204
$ws.Range("O8").Value = 0.3365562672414605
$ws.Range("P8").Value = 0.3365562672414606
$ws.Range("Q8").Value = 1208.418573813465
$ws.Range("R8").Value = 10875.76716432119
$ws.Range("S8").Value = 0.2529214954814828
$ws.Range("T8").Value = 0.2529214954814829

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 37.833149
$ws.Range("H9").Value = 113.499447
$ws.Range("I9").Value = 0.751498397443378
$ws.Range("J9").Value = 0.7514983974433781
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 38.58528233333333
$ws.Range("N9").Value = 115.755847
$ws.Range("O9").Value = 0.406569189096231
$ws.Range("P9").Value = 0.406569189096231
$ws.Range("Q9").Value = 1459.802735724068
$ws.Range("R9").Value = 13138.22462151661
$ws.Range("S9").Value = 0.3055360940556713
$ws.Range("T9").Value = 0.3055360940556713

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 37.833149
$ws.Range("H10").Value = 113.499447
$ws.Range("I10").Value = 0.751498397443378
$ws.Range("J10").Value = 0.7514983974433781
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 24.37857333333333
$ws.Range("N10").Value = 73.13571999999999
$ws.Range("O10").Value = 0.2568745436623085
$ws.Range("P10").Value = 0.2568745436623085
$ws.Range("Q10").Value = 922.3181973274266
$ws.Range("R10").Value = 8300.86377594684
$ws.Range("S10").Value = 0.1930408079062239
$ws.Range("T10").Value = 0.1930408079062239

